# Ran code for averaged intensities on spiral schemes.
#
# The "GossF" averaged-intensity table (HKL row labels in column B, ratios
# in C:M) gains three new sampling-scheme rows - "Spiral-90deg-10rot-5space",
# "Spiral-90deg-15rot-5space" and "Spiral-90deg-10rot-3space" - inserted
# right after the existing "Gaussian-Quadrature" row. Every row from
# "Gaussian-Quadrature" (previously the last row, r=16) through the end of
# the table shifts down accordingly, ending on row 19. Rows 1-9 (header +
# the first seven schemes) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 17:19 are brand-new sheet rows; give column A the same bold/bordered
# style ("s=1" in the original file) used by every other index cell in
# column A, by copying formats from an existing row of that kind.
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:A19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 10: index 8 - "Gaussian-Quadrature"
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.060443440083739
$ws.Range("D10").Value = 1.234423402731494
$ws.Range("E10").Value = 1.000829744036772
$ws.Range("F10").Value = 1.060443440083739
$ws.Range("G10").Value = 0.8587410004033036
$ws.Range("H10").Value = 1.451569371944928
$ws.Range("I10").Value = 0.962324404813809
$ws.Range("J10").Value = 1.234423402731494
$ws.Range("K10").Value = 1.117626573384134
$ws.Range("L10").Value = 1.089035006733936
$ws.Range("M10").Value = 1.094721894002341

# Row 11: index 9 - "Spiral-90deg-10rot-5space"
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 1.115807572074335
$ws.Range("D11").Value = 2.463641727999528
$ws.Range("E11").Value = 0.4348189935736123
$ws.Range("F11").Value = 1.115807572074335
$ws.Range("G11").Value = 1.678631161569662
$ws.Range("H11").Value = 0.07819509279137359
$ws.Range("I11").Value = 0.617226567013638
$ws.Range("J11").Value = 2.463641727999528
$ws.Range("K11").Value = 1.44923036078657
$ws.Range("L11").Value = 1.282518966430453
$ws.Range("M11").Value = 1.064720185837025

# Row 12: index 10 - "Spiral-90deg-15rot-5space"
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 1.10941712867094
$ws.Range("D12").Value = 2.469229007944313
$ws.Range("E12").Value = 0.434856367215113
$ws.Range("F12").Value = 1.10941712867094
$ws.Range("G12").Value = 1.682200625381976
$ws.Range("H12").Value = 0.07778761664708411
$ws.Range("I12").Value = 0.6163449794411291
$ws.Range("J12").Value = 2.469229007944313
$ws.Range("K12").Value = 1.452042687579713
$ws.Range("L12").Value = 1.280729908125326
$ws.Range("M12").Value = 1.064972620883426

# Row 13: index 11 - "Spiral-90deg-10rot-3space"
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 1.114252735605906
$ws.Range("D13").Value = 2.463901420304297
$ws.Range("E13").Value = 0.4349534686093818
$ws.Range("F13").Value = 1.114252735605906
$ws.Range("G13").Value = 1.679481002735638
$ws.Range("H13").Value = 0.07810873119124007
$ws.Range("I13").Value = 0.6172841120249255
$ws.Range("J13").Value = 2.463901420304297
$ws.Range("K13").Value = 1.449427444456839
$ws.Range("L13").Value = 1.281840090031373
$ws.Range("M13").Value = 1.064663578411898

# Row 14: index 12 - "NoRotation-tilt60deg"
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 2.213771999999995
$ws.Range("D14").Value = 1.523624000000005
$ws.Range("E14").Value = 0.4189959999999995
$ws.Range("F14").Value = 2.213771999999995
$ws.Range("G14").Value = 1.099428000000001
$ws.Range("H14").Value = 0.140872
$ws.Range("I14").Value = 0.7519960000000013
$ws.Range("J14").Value = 1.523624000000005
$ws.Range("K14").Value = 0.9713100000000023
$ws.Range("L14").Value = 1.592540999999999
$ws.Range("M14").Value = 1.024781333333334

# Row 15: index 13 - "Rotation-NoTilt"
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 3.434450000000003
$ws.Range("D15").Value = 0.01
$ws.Range("E15").Value = 0.45
$ws.Range("F15").Value = 3.434450000000003
$ws.Range("G15").Value = 0.3880625
$ws.Range("H15").Value = 0.22
$ws.Range("I15").Value = 0.9838875000000004
$ws.Range("J15").Value = 0.01
$ws.Range("K15").Value = 0.23
$ws.Range("L15").Value = 1.832225000000002
$ws.Range("M15").Value = 0.9144000000000005

# Row 16: index 14 - "Rotation-60detTilt"
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 2.385447859916797
$ws.Range("D16").Value = 0.4233949770752028
$ws.Range("E16").Value = 0.6843422908416021
$ws.Range("F16").Value = 2.385447859916797
$ws.Range("G16").Value = 0.6451964120064015
$ws.Range("H16").Value = 0.5584868202496011
$ws.Range("I16").Value = 0.9979515110399984
$ws.Range("J16").Value = 0.4233949770752028
$ws.Range("K16").Value = 0.5538686339584025
$ws.Range("L16").Value = 1.4696582469376
$ws.Range("M16").Value = 0.9491366451882671

# Row 17: index 15 - "HexGrid-90degTilt5degRes"
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9941439861373853
$ws.Range("D17").Value = 0.9882924057559362
$ws.Range("E17").Value = 0.9911795276129675
$ws.Range("F17").Value = 0.9941439861373853
$ws.Range("G17").Value = 0.9895893094430468
$ws.Range("H17").Value = 0.9872412484742678
$ws.Range("I17").Value = 0.9937465520539611
$ws.Range("J17").Value = 0.9882924057559362
$ws.Range("K17").Value = 0.9897359666844519
$ws.Range("L17").Value = 0.9919399764109186
$ws.Range("M17").Value = 0.9906988382462608

# Row 18: index 16 - "HexGrid-90degTilt22p5degRes"
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.9778878214107
$ws.Range("D18").Value = 1.081470585621645
$ws.Range("E18").Value = 0.9693686220836593
$ws.Range("F18").Value = 0.9778878214107
$ws.Range("G18").Value = 1.004877268752715
$ws.Range("H18").Value = 1.045522839282462
$ws.Range("I18").Value = 0.9714526093934206
$ws.Range("J18").Value = 1.081470585621645
$ws.Range("K18").Value = 1.025419603852652
$ws.Range("L18").Value = 1.001653712631676
$ws.Range("M18").Value = 1.008429957757433

# Row 19: index 17 - "HexGrid-60degTilt5degRes"
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9714474568022929
$ws.Range("D19").Value = 1.130027053825561
$ws.Range("E19").Value = 0.9669827382172724
$ws.Range("F19").Value = 0.9714474568022929
$ws.Range("G19").Value = 1.040957142282352
$ws.Range("H19").Value = 0.9786157634602403
$ws.Range("I19").Value = 0.9665210502646561
$ws.Range("J19").Value = 1.130027053825561
$ws.Range("K19").Value = 1.048504896021416
$ws.Range("L19").Value = 1.009976176411855
$ws.Range("M19").Value = 1.009091867475396
